# Applies the crypto-price/volume refresh described in the commit.
#
# Column D holds prices stored as literal text in the workbook (t="inlineStr"),
# even when the text happens to look like a plain number (e.g. "237.63" or
# "15.60"). Assigning such a string straight to Range.Value lets Excel
# auto-coerce it to a real number (dropping trailing zeros / switching to
# scientific notation), which would not match the source data. To keep the
# exact text we build the value with a text formula (="...") and then freeze
# it back to a plain value with Copy + PasteSpecial (xlPasteValues); this
# avoids touching the cell's NumberFormat/style (no new style gets created).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $escaped = $value -replace '"', '""'
    $rng.Formula = '="' + $escaped + '"'
    $rng.Copy()
    $rng.PasteSpecial(-4163)
}


# Row 2
$ws.Range("D2").Value = '26.168.00'
$ws.Range("E2").Value = '  +0.75%  '

# Row 3
$ws.Range("D3").Value = '1.752.25'
$ws.Range("E3").Value = '  +0.23%  '

# Row 4
Set-TextValue "D4" '0.9985'
$ws.Range("E4").Value = '  -0.27%  '

# Row 5
Set-TextValue "D5" '237.63'
$ws.Range("E5").Value = '  +1.17%  '

# Row 6
$ws.Range("B6").Value = 'XRP'
$ws.Range("C6").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextValue "D6" '0.5522'
$ws.Range("E6").Value = '  +6.03%  '

# Row 7
$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextValue "D7" '0.9991'
$ws.Range("E7").Value = '  -0.15%  '

# Row 8
Set-TextValue "D8" '0.2837'
$ws.Range("E8").Value = '  +0.17%  '

# Row 9
Set-TextValue "D9" '0.06192'
$ws.Range("E9").Value = '  +0.75%  '

# Row 10
$ws.Range("D10").Value = '1.748.04'
$ws.Range("E10").Value = '  -0.07%  '

# Row 11
Set-TextValue "D11" '0.07215'
$ws.Range("E11").Value = '  +2.64%  '

# Row 12
Set-TextValue "D12" '15.60'
$ws.Range("E12").Value = '  +0.78%  '

# Row 13
Set-TextValue "D13" '0.6536'
$ws.Range("E13").Value = '  +1.26%  '

# Row 14
Set-TextValue "D14" '4.652'
$ws.Range("E14").Value = '  +2.53%  '

# Row 15
Set-TextValue "D15" '78.69'
$ws.Range("E15").Value = '  +1.42%  '

# Row 16
Set-TextValue "D16" '0.9992'
$ws.Range("E16").Value = '  -0.09%  '

# Row 17
Set-TextValue "D17" '0.9979'
$ws.Range("E17").Value = '  -0.27%  '

# Row 18
$ws.Range("D18").Value = '26.044.45'

# Row 19
Set-TextValue "D19" '11.81'
$ws.Range("E19").Value = '  +2.62%  '

# Row 20
Set-TextValue "D20" '0.000006797'
$ws.Range("E20").Value = '  +2.54%  '

# Row 21
$ws.Range("D21").Value = '1.970.88'
$ws.Range("E21").Value = '  -0.36%  '

# Row 22
Set-TextValue "D22" '4.351'
$ws.Range("E22").Value = '  +4.42%  '

# Row 23
Set-TextValue "D23" '8.761'
$ws.Range("E23").Value = '  +1.20%  '

# Row 24
Set-TextValue "D24" '5.270'
$ws.Range("E24").Value = '  +2.13%  '

# Row 25
Set-TextValue "D25" '139.63'
$ws.Range("E25").Value = '  +0.35%  '

# Row 26
Set-TextValue "D26" '1.521'
$ws.Range("E26").Value = '  +0.95%  '

# Row 27
Set-TextValue "D27" '15.36'
$ws.Range("E27").Value = '  +1.69%  '

# Row 28
Set-TextValue "D28" '1.816'
$ws.Range("E28").Value = '  -1.32%  '

# Row 29
Set-TextValue "D29" '105.84'
$ws.Range("E29").Value = '  +2.88%  '

# Row 30
Set-TextValue "D30" '0.08421'
$ws.Range("E30").Value = '  +1.34%  '

# Row 31
Set-TextValue "D31" '3.818'
$ws.Range("E31").Value = '  +4.00%  '

# Row 32
Set-TextValue "D32" '3.661'
$ws.Range("E32").Value = '  +6.38%  '

# Row 33
Set-TextValue "D33" '0.04658'
$ws.Range("E33").Value = '  +4.02%  '

# Row 34
Set-TextValue "D34" '2.644'
$ws.Range("E34").Value = '  +1.13%  '

# Row 35
Set-TextValue "D35" '1.014'
$ws.Range("E35").Value = '  +2.67%  '

# Row 36
Set-TextValue "D36" '0.6352'
$ws.Range("E36").Value = '  +3.74%  '

# Row 37
$ws.Range("E37").Value = '  +1.13%  '

# Row 38
Set-TextValue "D38" '0.01631'
$ws.Range("E38").Value = '  +2.59%  '

# Row 39
Set-TextValue "D39" '1.980'
$ws.Range("E39").Value = '  +2.29%  '

# Row 40
Set-TextValue "D40" '0.9990'
$ws.Range("E40").Value = '  -0.08%  '

# Row 41
Set-TextValue "D41" '102.52'
$ws.Range("E41").Value = '  +1.50%  '

# Row 42
Set-TextValue "D42" '0.3973'
$ws.Range("E42").Value = '  +2.59%  '

# Row 43
Set-TextValue "D43" '0.7489'
$ws.Range("E43").Value = '  +1.88%  '

# Row 44
Set-TextValue "D44" '5.115'
$ws.Range("E44").Value = '  +0.68%  '

# Row 45
Set-TextValue "D45" '0.1158'
$ws.Range("E45").Value = '  +2.66%  '

# Row 46
Set-TextValue "D46" '6.409'
$ws.Range("E46").Value = '  +1.41%  '

# Row 47
Set-TextValue "D47" '0.05344'
$ws.Range("E47").Value = '  -2.31%  '

# Row 48
Set-TextValue "D48" '54.85'
$ws.Range("E48").Value = '  +3.40%  '

# Row 49
Set-TextValue "D49" '31.10'
$ws.Range("E49").Value = '  +3.47%  '

# Row 50
Set-TextValue "D50" '0.3517'
$ws.Range("E50").Value = '  +2.29%  '

# Row 51
Set-TextValue "D51" '7.641'
$ws.Range("E51").Value = '  +0.07%  '
